$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.251.91'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.348.68'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.24'
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("E9").Value = '  +2.96%  '
$ws.Range("E10").Value = '  +0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.05'
$ws.Range("E11").Value = '  +5.75%  '
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '689.23'
$ws.Range("E13").Value = '  +3.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.893.69'
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.311.89'
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.369.43'
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.21'
$ws.Range("E20").Value = '  +2.25%  '
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.50'
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '32.97'
$ws.Range("E28").Value = '  -1.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.53'
$ws.Range("E29").Value = '  +0.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.96'
$ws.Range("E30").Value = '  -5.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.07'
$ws.Range("E31").Value = '  +0.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '555.98'
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("E33").Value = '  +0.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '58.03'
$ws.Range("E34").Value = '  +2.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.708.19'
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.34'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.138'
$ws.Range("E38").Value = '  +4.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.82'
$ws.Range("E39").Value = '  +0.88%  '
$ws.Range("E40").Value = '  +1.72%  '
$ws.Range("E41").Value = '  -0.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₃0673'
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.24'
$ws.Range("E44").Value = '  -1.61%  '
$ws.Range("E45").Value = '  +1.29%  '
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("E49").Value = '  -1.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.99'
$ws.Range("E50").Value = '  +2.90%  '
$ws.Range("E51").Value = '  -2.17%  '
